# Remove the "Type" column (J) from the account-entries import sheet.
# The "Rule For" column (K) shifts left into J, taking over its comment-free
# state, formatting, and data; the "Type" column's dropdown data validation
# (Pool/All) and the associated header comment are removed along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting the entire column shifts K ("Rule For") and its data/styles left
# into J, and drops the data validation list that lived entirely in J2:J1025.
$ws.Range("J1").EntireColumn.Delete()

# The cell comment on the old J1 ("Type") is anchored to the absolute cell
# position, so it is NOT carried away by the column delete above - it keeps
# sitting on the now-shifted J1 ("Rule For"). Remove it explicitly, matching
# the target workbook where the new J1 has no comment (the old K1 never had
# one either).
$ws.Range("J1").Comment.Delete()

# Reflect the author's final selection resting on the new last column.
$ws.Range("J1:J1048576").Select()
